$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V2").Value = 1.54
$ws.Range("J3").Value = 2.37
$ws.Range("K4").Value = 1.92
$ws.Range("L4").Value = 2.87
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48
$ws.Range("G6").Value = 5.5
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 1.73
$ws.Range("L6").Value = 2.4
$ws.Range("R6").Value = 1.57
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("X6").Value = 26
$ws.Range("Y6").Value = 19
$ws.Range("AE6").Value = 21
$ws.Range("AI6").Value = 7
$ws.Range("AQ6").Value = 126
$ws.Range("AU6").Value = 9.5
$ws.Range("AW6").Value = 3.5
$ws.Range("AX6").Value = 9.5
$ws.Range("BA6").Value = 51
$ws.Range("K7").Value = 1.92
$ws.Range("R7").Value = 1.5
$ws.Range("AT7").Value = 2.37
$ws.Range("R9").Value = 1.44
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 1.57
$ws.Range("Q13").Value = 1.92
$ws.Range("R13").Value = 1.82
